# Generate Report for Handoff
# Updates the localization status workbook to reflect a new handoff run:
#  - the source file's GUID changes from 34d300fa-... to 87789cc7-...
#  - the generated .xlf files get a new content hash (6f164d87...)
#  - the "Latest Handoff" timestamps move forward a little

$wb = $excel.ActiveWorkbook

$oldGuid = "34d300fa-cdcc-4ded-b4ce-09f9691f9471"
$newGuid = "87789cc7-7006-41d3-a2e5-f7074c2323cd"

$oldHash = "800f5946412c46450b45ffc4cb89c0de01b119a6"
$newHash = "6f164d879a29e13a89d5058daa74b7cfabb6800c"

$newMdName    = "$newGuid.md"
$newZhCnXlf   = "$newGuid.$newHash.zh-cn.xlf"
$newDeDeXlf   = "$newGuid.$newHash.de-de.xlf"

$newOverviewDate = "2016-03-23 09:56:12"
$newZhCnDate      = "2016-03-23 09:56:03"

# NOTE: the hyperlink *targets* (the relationship Address) are left exactly
# as they were before the edit - only the cell text / displayed hyperlink
# text is refreshed to the new handoff file names.
$mdAddress     = "https://github.com/OpenLocalizationTest/oltest/blob/12bd272ec062cc19bde038ae193503a86dcbc669/e2e/$oldGuid.md"
$zhCnXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/46cab81d20817f3b3f8df807b9305defeb040d91/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/$oldGuid.$oldHash.zh-cn.xlf"
$deDeXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c31db7a0982dee35433435d007af135c1fa940b1/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/$oldGuid.$oldHash.de-de.xlf"

# ---------------------------------------------------------------------------
# Sheet "Overview": File Name (A2) + Latest Handoff Date (D2)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$a2 = $wsOverview.Range("A2")
$a2.Hyperlinks.Delete()
$a2.Value2 = $newMdName
$wsOverview.Hyperlinks.Add($a2, $mdAddress, "", "", $newMdName) | Out-Null

$wsOverview.Range("D2").Value2 = $newOverviewDate

# ---------------------------------------------------------------------------
# Sheet "zh-cn": File Name (A2), Latest Handoff File (D2), Latest Handoff
# Datetime (E2)
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$a2zh = $wsZhCn.Range("A2")
$a2zh.Hyperlinks.Delete()
$a2zh.Value2 = $newMdName
$wsZhCn.Hyperlinks.Add($a2zh, $mdAddress, "", "", $newMdName) | Out-Null

$d2zh = $wsZhCn.Range("D2")
$d2zh.Hyperlinks.Delete()
$d2zh.Value2 = $newZhCnXlf
$wsZhCn.Hyperlinks.Add($d2zh, $zhCnXlfAddress, "", "", $newZhCnXlf) | Out-Null

$wsZhCn.Range("E2").Value2 = $newZhCnDate

# ---------------------------------------------------------------------------
# Sheet "de-de": File Name (A2), Latest Handoff File (D2), Latest Handoff
# Datetime (E2, shares the same timestamp string as Overview!D2)
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$a2de = $wsDeDe.Range("A2")
$a2de.Hyperlinks.Delete()
$a2de.Value2 = $newMdName
$wsDeDe.Hyperlinks.Add($a2de, $mdAddress, "", "", $newMdName) | Out-Null

$d2de = $wsDeDe.Range("D2")
$d2de.Hyperlinks.Delete()
$d2de.Value2 = $newDeDeXlf
$wsDeDe.Hyperlinks.Add($d2de, $deDeXlfAddress, "", "", $newDeDeXlf) | Out-Null

$wsDeDe.Range("E2").Value2 = $newOverviewDate
